# "Small edit done for NodeURL input"
#
# The "Requirement" sheet lists upcoming Honda bike launches in a
# row-per-attribute / column-per-model layout. The "Honda Grazia" column
# (column A: model name, price, expected-launch date) is removed, shifting
# the remaining three models (Honda CB Hornet 160R, Honda XBlade,
# Honda CBF190R) one column to the left. The price text for what is now
# the last column is also corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Honda Grazia" column entirely; B:D (Hornet/XBlade/CBF190R)
# shift left into A:C.
$ws.Columns("A").Delete()

# "1.10 lakh onwards" -> "1.1 Lakh onwards" for the CBF190R price, which is
# now in column C, row 2.
$ws.Range("C2").Value = "1.1 Lakh onwards"
